# Apply corrections to the loan-schedule workbook:
#  - Summary sheet: update the active selection/cursor position
#  - Repayment schedule sheet: remove the now-unused "Over Due" (column O)
#    values for the amortization rows
#  - Transactions sheet: update the active selection/cursor position and
#    correct several transaction figures

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - just move the saved selection from A7:XFD13 to B5
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Select() | Out-Null

# ---------------------------------------------------------------------
# Repayment schedule sheet - clear out column O (rows 2-14), fully
# removing the cells (contents + formatting) rather than just blanking
# their values, leaving column P untouched.
# ---------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$oRange = $wsSchedule.Range("O2:O14")
$oRange.ClearFormats()
$oRange.Value = $null

# ---------------------------------------------------------------------
# Transactions sheet - correct several figures and move the selection
# ---------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")

$wsTransactions.Range("A2").Value = 36
$wsTransactions.Range("J2").Value = 9133.2199999999993

$wsTransactions.Range("A3").Value = 34
$wsTransactions.Range("C3").Value = 42064
$wsTransactions.Range("E3").Value = 963.77
$wsTransactions.Range("F3").Value = 866.78
$wsTransactions.Range("G3").Value = 96.99
$wsTransactions.Range("J3").Value = 4133.22

$wsTransactions.Range("A4").Value = 32

$wsTransactions.Range("A2").Select() | Out-Null
